# chore: 新增不需要校验 token API
# Populate the "已有用户账号" (existing user accounts) list in column A
# with the updated set of sample user accounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full, final ordered list of user account names that should occupy
# A2:A26 after the edit (existing accounts plus newly added ones,
# interleaved in the same order the commit introduced them).
$accounts = @(
    "qishilong",
    "教授",
    "王师傅",
    "田师傅",
    "老谈",
    "22级学弟",
    "haha",
    "C++",
    "php",
    "students",
    "xxxxx",
    "xxxxfadfa",
    "nihaonihao",
    "666",
    "原神",
    "期末60",
    "chihuo",
    "大吃货",
    "小吃货",
    "我的中路特别稳",
    "我的上路特别稳",
    "wanghaha",
    "nihao",
    "111",
    "nihao11"
)

for ($i = 0; $i -lt $accounts.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $text = $accounts[$i]

    $looksNumeric = $text -match '^[0-9]+$'

    if ($looksNumeric) {
        # Force Excel to store purely-numeric-looking account names (e.g. "666")
        # as text rather than silently converting them to a Number cell, then
        # drop the temporary text format so no stray formatting lingers.
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}
